# Update the build/version string throughout the workbook for the
# "Coal Mine Boundaries and Methane Sources - version 1.0.0" release.
#
# Old version tag: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version tag: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: " + $newVersion

$about.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for No. 7 Coal Mine, United States, M1051, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---------------------------
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S holds the build_version for every data row (rows 2-33).
$lastRow = $data.Cells.Item($data.Rows.Count, 19).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $data.Cells.Item($r, 19).Value = $newVersion
}
